$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replicate the original authoring order so new shared strings are
# interned in the same sequence as the authored workbook.
$ws.Range("H1").Value = "Labor Booking User"
$ws.Range("H2").Value = "a811K0000004fpN"

$ws.Range("I1").Value = "SiteID"
$ws.Range("J1").Value = "Location ID"
$ws.Range("K1").Value = "Location Number"

$ws.Range("I2").Value = "a7q410000004I1W"
$ws.Range("J2").Value = "a7Z4100000000hb"
$ws.Range("K2").Value = "SY_ReceiptLoc"

# Row 3 duplicates row 2's values
$ws.Range("H3").Value = "a811K0000004fpN"
$ws.Range("I3").Value = "a7q410000004I1W"
$ws.Range("J3").Value = "a7Z4100000000hb"
$ws.Range("K3").Value = "SY_ReceiptLoc"

# Resize the columns touched by this change to their authored best-fit
# widths. (ColumnWidth is expressed in characters of the Normal font.)
$ws.Columns("B").ColumnWidth = 20.666666666666668
$ws.Columns("H").ColumnWidth = 16.166666666666668
$ws.Columns("I").ColumnWidth = 15.666666666666666
$ws.Columns("J").ColumnWidth = 15.166666666666666
$ws.Columns("K").ColumnWidth = 14.333333333333334

# Update the active selection to match the authored state
$ws.Range("J4").Select() | Out-Null
